# Updates vm_pu.xlsx results for "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.060590024698926
$ws.Range("D2").Value = 1.060778303270036
$ws.Range("E2").Value = 1.074146206274934
$ws.Range("F2").Value = 1.081241313296882
$ws.Range("I2").Value = 1.052663900367918
$ws.Range("J2").Value = 1.065570088629278
$ws.Range("K2").Value = 1.063504157313564
$ws.Range("L2").Value = 1.07683619141211
$ws.Range("M2").Value = 1.083912649957102
$ws.Range("N2").Value = 1.025320154627256

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.061888174758018
$ws.Range("D3").Value = 1.061808859511142
$ws.Range("E3").Value = 1.075470313933804
$ws.Range("F3").Value = 1.082707631296825
$ws.Range("I3").Value = 1.053125302292105
$ws.Range("J3").Value = 1.066520283334186
$ws.Range("K3").Value = 1.064348660370973
$ws.Range("L3").Value = 1.077976123498642
$ws.Range("M3").Value = 1.085195807296432
$ws.Range("N3").Value = 1.02566169062054

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.062727169861754
$ws.Range("D4").Value = 1.062474763606663
$ws.Range("E4").Value = 1.076326481620599
$ws.Range("F4").Value = 1.083656074218143
$ws.Range("I4").Value = 1.05342202551552
$ws.Range("J4").Value = 1.067133619399415
$ws.Range("K4").Value = 1.064893569850249
$ws.Range("L4").Value = 1.078712567576565
$ws.Range("M4").Value = 1.086025209102087
$ws.Range("N4").Value = 1.025881605000359

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.063079649562544
$ws.Range("D5").Value = 1.062754488660868
$ws.Range("E5").Value = 1.076686270458313
$ws.Range("F5").Value = 1.084054717347069
$ws.Range("I5").Value = 1.053546330447778
$ws.Range("J5").Value = 1.067391109072801
$ws.Range("K5").Value = 1.065122283669155
$ws.Range("L5").Value = 1.079021892652779
$ws.Range("M5").Value = 1.086373681437482
$ws.Range("N5").Value = 1.02597379893786

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.06313881886788
$ws.Range("D6").Value = 1.062801442846632
$ws.Range("E6").Value = 1.07674667225613
$ws.Range("F6").Value = 1.084121646591251
$ws.Range("I6").Value = 1.053567176186708
$ws.Range("J6").Value = 1.067434321889016
$ws.Range("K6").Value = 1.065160664333181
$ws.Range("L6").Value = 1.079073813614629
$ws.Range("M6").Value = 1.086432179387368
$ws.Range("N6").Value = 1.025989263602133

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.062731880626229
$ws.Range("D7").Value = 1.062478502171609
$ws.Range("E7").Value = 1.076331289697707
$ws.Range("F7").Value = 1.083661401224699
$ws.Range("I7").Value = 1.053423688201102
$ws.Range("J7").Value = 1.067137061387168
$ws.Range("K7").Value = 1.064896627370808
$ws.Range("L7").Value = 1.07871670187337
$ws.Range("M7").Value = 1.086029866219383
$ws.Range("N7").Value = 1.025882837913012

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.061028949167522
$ws.Range("D8").Value = 1.061126779767938
$ws.Range("E8").Value = 1.074593825144412
$ws.Range("F8").Value = 1.08173694023728
$ws.Range("I8").Value = 1.052820214223725
$ws.Range("J8").Value = 1.065891524558308
$ws.Range("K8").Value = 1.063789881647677
$ws.Range("L8").Value = 1.077221680887792
$ws.Range("M8").Value = 1.084346485086298
$ws.Range("N8").Value = 1.025435803018751

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.058020353852263
$ws.Range("D9").Value = 1.058737590022768
$ws.Range("E9").Value = 1.07152727331446
$ws.Range("F9").Value = 1.078342808600583
$ws.Range("I9").Value = 1.051742687460914
$ws.Range("J9").Value = 1.063685074755951
$ws.Range("K9").Value = 1.061827736259041
$ws.Range("L9").Value = 1.074578151596086
$ws.Range("M9").Value = 1.081373175859852
$ws.Range("N9").Value = 1.02463973565859

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.056009098404024
$ws.Range("D10").Value = 1.057139725858444
$ws.Range("E10").Value = 1.069479345261979
$ws.Range("F10").Value = 1.076077756534028
$ws.Range("I10").Value = 1.051014728012661
$ws.Range("J10").Value = 1.062206074426997
$ws.Range("K10").Value = 1.060511456873412
$ws.Range("L10").Value = 1.072809441293093
$ws.Range("M10").Value = 1.079386018772014
$ws.Range("N10").Value = 1.024103354979277

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.055136830860946
$ws.Range("D11").Value = 1.056446590642913
$ws.Range("E11").Value = 1.06859166822972
$ws.Range("F11").Value = 1.075096353818934
$ws.Range("I11").Value = 1.05069720990746
$ws.Range("J11").Value = 1.061563700873484
$ws.Range("K11").Value = 1.059939514308247
$ws.Range("L11").Value = 1.072042010200772
$ws.Range("M11").Value = 1.078524325020247
$ws.Range("N11").Value = 1.023869736770973

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.054812619069554
$ws.Range("D12").Value = 1.056188938135335
$ws.Range("E12").Value = 1.068261803628373
$ws.Range("F12").Value = 1.074731718273999
$ws.Range("I12").Value = 1.050578920772237
$ws.Range("J12").Value = 1.061324796843227
$ws.Range("K12").Value = 1.059726767779517
$ws.Range("L12").Value = 1.071756711994314
$ws.Range("M12").Value = 1.078204061929647
$ws.Range("N12").Value = 1.023782754642839

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.054882173314834
$ws.Range("D13").Value = 1.056244214182764
$ws.Range("E13").Value = 1.068332567203845
$ws.Range("F13").Value = 1.074809938402489
$ws.Range("I13").Value = 1.050604310011371
$ws.Range("J13").Value = 1.061376056112149
$ws.Range("K13").Value = 1.059772416334449
$ws.Range("L13").Value = 1.071817920368137
$ws.Range("M13").Value = 1.078272768205527
$ws.Range("N13").Value = 1.023801421944396

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.055110035785042
$ws.Range("D14").Value = 1.056425296927581
$ws.Range("E14").Value = 1.068564404441185
$ws.Range("F14").Value = 1.075066214981216
$ws.Range("I14").Value = 1.050687439218246
$ws.Range("J14").Value = 1.061543959091482
$ws.Range("K14").Value = 1.059921934785895
$ws.Range("L14").Value = 1.072018432287136
$ws.Range("M14").Value = 1.078497855912902
$ws.Range("N14").Value = 1.023862551009586

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.055250401016062
$ws.Range("D15").Value = 1.056536842538523
$ws.Range("E15").Value = 1.068707228070839
$ws.Range("F15").Value = 1.075224102166111
$ws.Range("I15").Value = 1.050738611584516
$ws.Range("J15").Value = 1.061647370067322
$ws.Range("K15").Value = 1.060014017972808
$ws.Range("L15").Value = 1.072141942341461
$ws.Range("M15").Value = 1.078636514327112
$ws.Range("N15").Value = 1.023900187317574

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.056066958382333
$ws.Range("D16").Value = 1.057185700335952
$ws.Range("E16").Value = 1.069538237849324
$ws.Range("F16").Value = 1.076142875446971
$ws.Range("I16").Value = 1.05103575185371
$ws.Range("J16").Value = 1.062248665085115
$ws.Range("K16").Value = 1.060549372709686
$ws.Range("L16").Value = 1.072860339744207
$ws.Range("M16").Value = 1.079443179952062
$ws.Range("N16").Value = 1.024118830645252

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.056578789827147
$ws.Range("D17").Value = 1.057592374698114
$ws.Range("E17").Value = 1.070059260986899
$ws.Range("F17").Value = 1.076719027586864
$ws.Range("I17").Value = 1.05122152093696
$ws.Range("J17").Value = 1.062625314947135
$ws.Range("K17").Value = 1.060884652792968
$ws.Range("L17").Value = 1.073310548471345
$ws.Range("M17").Value = 1.079948843959599
$ws.Range("N17").Value = 1.024255614293931

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.056877199599829
$ws.Range("D18").Value = 1.0578294606545
$ws.Range("E18").Value = 1.070363077442127
$ws.Range("F18").Value = 1.077055027927091
$ws.Range("I18").Value = 1.051329654470964
$ws.Range("J18").Value = 1.062844819917804
$ws.Range("K18").Value = 1.061080024708867
$ws.Range("L18").Value = 1.0735729966404
$ws.Range("M18").Value = 1.080243669770158
$ws.Range("N18").Value = 1.024335266574813

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.056978927248476
$ws.Range("D19").Value = 1.057910280535986
$ws.Range("E19").Value = 1.070466656306894
$ws.Range("F19").Value = 1.077169585446985
$ws.Range("I19").Value = 1.051366487549758
$ws.Range("J19").Value = 1.062919633581918
$ws.Range("K19").Value = 1.061146609140005
$ws.Range("L19").Value = 1.073662459215283
$ws.Range("M19").Value = 1.080344177646078
$ws.Range("N19").Value = 1.024362403712006

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05652388893857
$ws.Range("D20").Value = 1.057548754859414
$ws.Range("E20").Value = 1.070003369226071
$ws.Range("F20").Value = 1.076657218107465
$ws.Range("I20").Value = 1.051201612695893
$ws.Range("J20").Value = 1.062584923510749
$ws.Range("K20").Value = 1.060848700224049
$ws.Range("L20").Value = 1.073262260985069
$ws.Range("M20").Value = 1.079894603390952
$ws.Range("N20").Value = 1.024240952294683

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.055042941884062
$ws.Range("D21").Value = 1.056371977879655
$ws.Range("E21").Value = 1.068496138108687
$ws.Range("F21").Value = 1.074990750691776
$ws.Range("I21").Value = 1.050662969384295
$ws.Range("J21").Value = 1.06149452410382
$ws.Range("K21").Value = 1.059877913681352
$ws.Range("L21").Value = 1.071959393190221
$ws.Range("M21").Value = 1.078431578525408
$ws.Range("N21").Value = 1.023844555718654

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.054110577891894
$ws.Range("D22").Value = 1.055630983356001
$ws.Range("E22").Value = 1.067547659011839
$ws.Range("F22").Value = 1.073942401556186
$ws.Range("I22").Value = 1.050322283965019
$ws.Range("J22").Value = 1.060807221255273
$ws.Range("K22").Value = 1.059265795029911
$ws.Range("L22").Value = 1.071138836901956
$ws.Range("M22").Value = 1.077510605257875
$ws.Range("N22").Value = 1.023594133387215

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.054604960430724
$ws.Range("D23").Value = 1.056023904786418
$ws.Range("E23").Value = 1.068050545423278
$ws.Range("F23").Value = 1.07449820767419
$ws.Range("I23").Value = 1.050503079863899
$ws.Range("J23").Value = 1.061171738369132
$ws.Range("K23").Value = 1.059590457402811
$ws.Range("L23").Value = 1.071573962756176
$ws.Range("M23").Value = 1.077998937683806
$ws.Range("N23").Value = 1.023727000441111

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.056548696697113
$ws.Range("D24").Value = 1.057568465158543
$ws.Range("E24").Value = 1.070028624577134
$ws.Range("F24").Value = 1.076685147334037
$ws.Range("I24").Value = 1.051210609060243
$ws.Range("J24").Value = 1.062603175244511
$ws.Range("K24").Value = 1.060864946233247
$ws.Range("L24").Value = 1.073284080487813
$ws.Range("M24").Value = 1.079919112737961
$ws.Range("N24").Value = 1.024247577826694

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.058799101661476
$ws.Range("D25").Value = 1.05935613375236
$ws.Range("E25").Value = 1.072320659272103
$ws.Range("F25").Value = 1.079220656594819
$ws.Range("I25").Value = 1.052022939450397
$ws.Range("J25").Value = 1.064256896369458
$ws.Range("K25").Value = 1.062336427608576
$ws.Range("L25").Value = 1.075262670785215
$ws.Range("M25").Value = 1.082142701047517
$ws.Range("N25").Value = 1.024846532567328

